$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that change, so numeric-looking strings
# (e.g. thousands-dot formatted prices) are preserved exactly as text,
# matching the original inline-string cell content instead of being
# auto-converted to numbers by Excel's smart input parsing.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values from the crypto price refresh
$ws.Range('D2').Value = '29.755.31'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.607.03'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '213.21'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').Value = '28.29'
$ws.Range('E8').Value = '  +5.59%  '
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').Value = '0.0910'
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').Value = '1.835.87'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '1.597.53'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('E14').Value = '  +3.89%  '
$ws.Range('D15').Value = '29.739.37'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').Value = '64.17'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = '242.36'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').Value = '7.88'
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').Value = '4.04'
$ws.Range('D23').Value = '9.41'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').Value = '155.16'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').Value = '15.49'
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('E31').Value = '  +0.81%  '
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('D34').Value = '1.426.10'
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  +3.37%  '
$ws.Range('E36').Value = '  +2.37%  '
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('D40').Value = '0.548'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('D41').Value = '56.91'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '0.0495'
$ws.Range('E42').Value = '  +5.82%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.96'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '0.818'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '66.55'
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '0.984'
$ws.Range('E47').Value = '  +17.48%  '
$ws.Range('D49').Value = '1.745.24'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = '86.58'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('E51').Value = '  -0.07%  '
